$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (B4): append a note to the Animal activity text
$ws.Range("B4").Value = "Construir  Interfaces de cadastro e consulta Animal (Jaime e valter) falta actualizar dados"

# Row 3 (B3): remove "falta Deletar" from the Cliente activity text
$ws.Range("B3").Value = "Construir  Interfaces de cadastro e consulta Cliente  (Jaime) "

# Row 3 (C3): mark this activity as "ok"
$ws.Range("C3").Value = "ok"

# Update the active selection to C4
$ws.Range("C4").Select()
